$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4399.2
$ws.Range("I51").Value = 3998.6667
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 3998.6667
$ws.Range("L51").Value = 5000
$ws.Range("M51").Value = -3514.6667
$ws.Range("N51").Value = -5968

$ws.Range("H138").Value = 2987.873
$ws.Range("J138").Value = 3517.6191
$ws.Range("L138").Value = 10552.8573
$ws.Range("N138").Value = -20832.8573

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3214.7112
$ws.Range("I63").Value = 1336.4
$ws.Range("K63").Value = 1336.4
$ws.Range("M63").Value = -650.4000000000001

$ws.Range("H66").Value = 3214.7112
$ws.Range("I66").Value = 1336.4
$ws.Range("K66").Value = 6682
$ws.Range("M66").Value = -3250

$ws.Range("H74").Value = 3883.05
$ws.Range("I74").Value = 1189.2142
$ws.Range("K74").Value = 1189.2142
$ws.Range("M74").Value = -315.2141999999999

$ws.Range("H77").Value = 3883.05
$ws.Range("I77").Value = 1189.2142
$ws.Range("K77").Value = 5946.071
$ws.Range("M77").Value = -1578.071

$ws.Range("H132").Value = 1783.3684
$ws.Range("I132").Value = 1632.9375
$ws.Range("K132").Value = 4898.8125
$ws.Range("M132").Value = -2368.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2035.6154
$ws.Range("I86").Value = 1746.6
$ws.Range("K86").Value = 1746.6
$ws.Range("M86").Value = -623.5999999999999

$ws.Range("H89").Value = 2035.6154
$ws.Range("I89").Value = 1746.6
$ws.Range("K89").Value = 8733
$ws.Range("M89").Value = -3117

$ws.Range("H134").Value = 4306.2705
$ws.Range("I134").Value = 3879.8076
$ws.Range("K134").Value = 11639.4228
$ws.Range("M134").Value = -9104.4228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2246.353
$ws.Range("J16").Value = 3535.5
$ws.Range("L16").Value = 3535.5
$ws.Range("N16").Value = -4109.5

$ws.Range("H18").Value = 50000
$ws.Range("J18").Value = 50000
$ws.Range("L18").Value = 50000
$ws.Range("N18").Value = -50460

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H86").Value = 25759910
$ws.Range("I86").Value = 43591890
$ws.Range("J86").Value = 2608.6667
$ws.Range("K86").Value = 43591890
$ws.Range("L86").Value = 2608.6667
$ws.Range("M86").Value = -43590767
$ws.Range("N86").Value = -4854.6667

$ws.Range("H88").Value = 30558.572
$ws.Range("J88").Value = 30558.572
$ws.Range("L88").Value = 30558.572
$ws.Range("N88").Value = -31370.572

$ws.Range("H89").Value = 25759910
$ws.Range("I89").Value = 43591890
$ws.Range("J89").Value = 2608.6667
$ws.Range("K89").Value = 217959450
$ws.Range("L89").Value = 13043.3335
$ws.Range("M89").Value = -217953834
$ws.Range("N89").Value = -24275.3335

$ws.Range("H91").Value = 30558.572
$ws.Range("J91").Value = 30558.572
$ws.Range("L91").Value = 30558.572
$ws.Range("N91").Value = -33366.572

$ws.Range("H99").Value = 10019.216
$ws.Range("I99").Value = 5567.7617
$ws.Range("J99").Value = 15861.75
$ws.Range("K99").Value = 5567.7617
$ws.Range("L99").Value = 15861.75
$ws.Range("M99").Value = -4069.7617
$ws.Range("N99").Value = -18857.75

$ws.Range("H100").Value = 83725
$ws.Range("J100").Value = 83725
$ws.Range("L100").Value = 83725
$ws.Range("N100").Value = -85889

$ws.Range("H107").Value = 1007.26086
$ws.Range("I107").Value = 951.7368
$ws.Range("K107").Value = 951.7368
$ws.Range("M107").Value = 968.2632

$ws.Range("H108").Value = 66666.336
$ws.Range("J108").Value = 66666.336
$ws.Range("L108").Value = 66666.336
$ws.Range("N108").Value = -74346.336

$ws.Range("H109").Value = 99999
$ws.Range("J109").Value = 99999
$ws.Range("L109").Value = 99999
$ws.Range("N109").Value = -102079

$ws.Range("H110").Value = 99999.5
$ws.Range("J110").Value = 99999.5
$ws.Range("L110").Value = 99999.5
$ws.Range("N110").Value = -108179.5

$ws.Range("H111").Value = 199976.33
$ws.Range("J111").Value = 199976.33
$ws.Range("L111").Value = 199976.33
$ws.Range("N111").Value = -208156.33

$ws.Range("H112").Value = 83332.664
$ws.Range("J112").Value = 83332.664
$ws.Range("L112").Value = 83332.664
$ws.Range("N112").Value = -86286.664

$ws.Range("H113").Value = 2246.353
$ws.Range("J113").Value = 3535.5
$ws.Range("L113").Value = 3535.5
$ws.Range("N113").Value = -7875.5

$ws.Range("H114").Value = 72413
$ws.Range("J114").Value = 72413
$ws.Range("L114").Value = 72413
$ws.Range("N114").Value = -81091

$ws.Range("H116").Value = 95000
$ws.Range("J116").Value = 95000
$ws.Range("L116").Value = 95000
$ws.Range("N116").Value = -104178

$ws.Range("H121").Value = 65000
$ws.Range("J121").Value = 65000
$ws.Range("L121").Value = 65000
$ws.Range("N121").Value = -67620

$ws.Range("H124").Value = 47801
$ws.Range("J124").Value = 47801
$ws.Range("L124").Value = 47801
$ws.Range("N124").Value = -52711

$ws.Range("H126").Value = 10019.216
$ws.Range("I126").Value = 5567.7617
$ws.Range("J126").Value = 15861.75
$ws.Range("K126").Value = 16703.2851
$ws.Range("L126").Value = 47585.25
$ws.Range("M126").Value = -14233.2851
$ws.Range("N126").Value = -52525.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3657.4443
$ws.Range("I122").Value = 1654.5
$ws.Range("K122").Value = 4963.5
$ws.Range("M122").Value = -2513.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3778.1904
$ws.Range("I136").Value = 3509.5715
$ws.Range("K136").Value = 10528.7145
$ws.Range("M136").Value = -7978.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3230.3
$ws.Range("I126").Value = 3014.4285
$ws.Range("K126").Value = 9043.2855
$ws.Range("M126").Value = -6573.2855

$ws.Range("H132").Value = 2961.6445
$ws.Range("I132").Value = 2402.639
$ws.Range("J132").Value = 5197.6665
$ws.Range("K132").Value = 7207.917
$ws.Range("L132").Value = 15592.9995
$ws.Range("M132").Value = -4677.917
$ws.Range("N132").Value = -20652.9995
